$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2147
$ws.Range("B3").Value = 16589
$ws.Range("B4").Value = 5752
$ws.Range("B5").Value = 34.67
$ws.Range("B6").Value = 7.73
$ws.Range("B7").Value = 6.29
